$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (single "." as decimal separator)
# must be forced to Text format first, otherwise Excel auto-converts them to a
# number (dropping the exact printed text, e.g. trailing zeros: "1.000" -> 1).
$textCells = @(
    "D4", "D5", "D8", "D9", "D10", "D11", "D12", "D14",
    "D15", "D16", "D18", "D19", "D20", "D22", "D23", "D24",
    "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33",
    "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43",
    "D44", "D45", "D46", "D47", "D48", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cryptos snapshot values
$ws.Range("D2").Value = "30.441.73"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.916.57"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "243.71"
$ws.Range("E5").Value = "  +2.08%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("D8").Value = "0.2868"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").Value = "0.06841"
$ws.Range("E9").Value = "  +4.03%  "
$ws.Range("D10").Value = "110.54"
$ws.Range("E10").Value = "  +11.31%  "
$ws.Range("D11").Value = "18.43"
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").Value = "0.07727"
$ws.Range("E12").Value = "  +1.90%  "
$ws.Range("D13").Value = "1.890.13"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "5.286"
$ws.Range("E14").Value = "  +3.20%  "
$ws.Range("D15").Value = "0.6569"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "296.88"
$ws.Range("E16").Value = "  -2.93%  "
$ws.Range("D17").Value = "30.439.41"
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "0.000007633"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "12.96"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").Value = "2.141.28"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "1.0000"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "5.242"
$ws.Range("E23").Value = "  +2.37%  "
$ws.Range("D24").Value = "6.212"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Value = "9.374"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("D26").Value = "21.77"
$ws.Range("E26").Value = "  +6.83%  "
$ws.Range("D27").Value = "168.86"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").Value = "2.096"
$ws.Range("E28").Value = "  +7.76%  "
$ws.Range("D29").Value = "0.1069"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").Value = "1.365"
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").Value = "3.981"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "0.05035"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("D35").Value = "0.7350"
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").Value = "0.02074"
$ws.Range("E36").Value = "  +6.37%  "
$ws.Range("D37").Value = "2.740"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").Value = "2.677"
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").Value = "109.55"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("D41").Value = "0.8702"
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("D42").Value = "5.852"
$ws.Range("E42").Value = "  +4.08%  "
$ws.Range("D43").Value = "0.4251"
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "51.61"
$ws.Range("E45").Value = "  +20.53%  "
$ws.Range("D46").Value = "67.33"
$ws.Range("E46").Value = "  +2.10%  "
$ws.Range("D47").Value = "7.194"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").Value = "9.248"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("D50").Value = "34.97"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").Value = "0.2443"
$ws.Range("E51").Value = "  +10.70%  "
